# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing": insert a new row for D.Lock (right after T.Bridgewater),
# update a handful of stat lines, and renumber the index column.
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Insert new row 3 for D.Lock; pushes M.Gordon..D.Hamilton down one row.
$rushing.Rows("3").Insert()

# Re-apply the header-style formatting (bold, thin border, centered) that the
# rest of column A uses, matching the existing rows.
$rushing.Range("A3").Font.Bold = $true
$rushing.Range("A3").HorizontalAlignment = -4108
$rushing.Range("A3").VerticalAlignment = -4160
$rushing.Range("A3").Borders.LineStyle = 1

$rushing.Range("B3").Value = "D.Lock"
$rushing.Range("C3").Value = 0
$rushing.Range("D3").Value = 1
$rushing.Range("E3").Value = 0
$rushing.Range("F3").Value = 1

# Updated weekly stats for returning players.
$rushing.Range("D2").Value = 7
$rushing.Range("E2").Value = 13

$rushing.Range("C4").Value = 93
$rushing.Range("D4").Value = 67
$rushing.Range("E4").Value = 10
$rushing.Range("F4").Value = 34

$rushing.Range("C5").Value = 92
$rushing.Range("D5").Value = 57
$rushing.Range("E5").Value = 20
$rushing.Range("F5").Value = 21

# Renumber the index column (A) sequentially now that a row was inserted.
$rushing.Range("A2").Value = 0
$rushing.Range("A3").Value = 1
$rushing.Range("A4").Value = 2
$rushing.Range("A5").Value = 3
$rushing.Range("A6").Value = 4
$rushing.Range("A7").Value = 5
$rushing.Range("A8").Value = 6
$rushing.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# Sheet "Receiving": update weekly stats (no rows added/removed).
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 30
$receiving.Range("D2").Value = 21

$receiving.Range("C3").Value = 46
$receiving.Range("D3").Value = 36

$receiving.Range("C5").Value = 55
$receiving.Range("D5").Value = 42
$receiving.Range("E5").Value = 30

$receiving.Range("C6").Value = 72
$receiving.Range("E6").Value = 22

$receiving.Range("C7").Value = 54
$receiving.Range("D7").Value = 37
$receiving.Range("E7").Value = 17
$receiving.Range("F7").Value = 8

$receiving.Range("C9").Value = 5

$receiving.Range("C10").Value = 71
$receiving.Range("D10").Value = 57
$receiving.Range("E10").Value = 12
$receiving.Range("F10").Value = 7

$receiving.Range("C11").Value = 32
$receiving.Range("D11").Value = 29
$receiving.Range("E11").Value = 5
$receiving.Range("F11").Value = 3
